$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 659
    $ws.Range("F3").Value = 3914
    $ws.Range("F5").Value = 732
    $ws.Range("G5").Value = 35
}
